$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price + volume figures).
# Column D (Price) values must remain plain text exactly as scraped (they can
# contain thousand-separator dots or trailing zeros), so we briefly force a
# text number format before assigning, then restore the default style so the
# cell keeps its original (unstyled) appearance.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.231.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.912.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "363.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.539"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.99%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -6.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.26%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0834"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.370.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.915.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.955"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.136.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("E20").Value = "  -3.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0948"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.93%  "
$ws.Range("E27").Value = "  -5.75%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.79%  "
$ws.Range("E30").Value = "  -6.69%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.16%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "35.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.53%  "
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("E41").Value = "  -7.04%  "
$ws.Range("E42").Value = "  -6.64%  "
$ws.Range("E43").Value = "  -5.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "118.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.062.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("E48").Value = "  -7.72%  "
$ws.Range("E49").Value = "  -8.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.199.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  -6.09%  "
